# Auto-derived script applying the scraper-refresh diff to 广州-漫展信息.xlsx
# Sheet '展览' (Exhibitions, sheet1) and sheet '全部类型' (All types, sheet4) each
# lose their duplicate '第九届娃展沙龙' row, causing every later row's B:I content
# to shift up by one (column A keeps its original static index numbers), and a
# number of 'want to go' (F column) counts are refreshed to newer scraped values.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item('展览')
$ws1.Cells.Item(6, 6).Value = 492
$ws1.Cells.Item(7, 6).Value = 492
$ws1.Cells.Item(8, 6).Value = 766
$ws1.Cells.Item(9, 6).Value = 180
$ws1.Cells.Item(10, 6).Value = 1364
$ws1.Cells.Item(11, 6).Value = 793
$ws1.Cells.Item(13, 6).Value = 597
$ws1.Cells.Item(14, 6).Value = 139
$ws1.Cells.Item(16, 6).Value = 18
$ws1.Cells.Item(18, 6).Value = 95

# Shift rows 19-33 up by one logical event (content that used to occupy the next
# row down), refreshing F/G counts where the source re-scrape changed them.
$ws1.Cells.Item(19, 2).Value = '2024-03-30'
$ws1.Cells.Item(19, 3).Value = '广州·AP动漫游戏嘉年华'
$ws1.Cells.Item(19, 4).Value = '新港东路630-638号 南丰国际会展中心'
$ws1.Cells.Item(19, 5).Value = '2024.03.30 09:00-03.31 17:00'
$ws1.Cells.Item(19, 6).Value = 1470
$ws1.Cells.Item(19, 7).Value = '不可售'
$ws1.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82121'
$ws1.Cells.Item(19, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/SLJ05mdG1709019165508.jpeg'

$ws1.Cells.Item(20, 2).Value = '2024-04-04'
$ws1.Cells.Item(20, 3).Value = '广州·Look Look动漫嘉年华'
$ws1.Cells.Item(20, 4).Value = '东沙大道16号 健康方舟'
$ws1.Cells.Item(20, 5).Value = '2024.04.04 10:00-04.05 17:30'
$ws1.Cells.Item(20, 6).Value = 171
$ws1.Cells.Item(20, 7).Value = 29.9
$ws1.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82319'
$ws1.Cells.Item(20, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png'

$ws1.Cells.Item(21, 2).Value = '2024-04-05'
$ws1.Cells.Item(21, 3).Value = '广州·允月秘境动漫嘉年华'
$ws1.Cells.Item(21, 4).Value = '人民北路686号广东广播中心大楼 广东广播电视台(人民北路)'
$ws1.Cells.Item(21, 5).Value = '2024.04.05 11:30-04.06 18:00'
$ws1.Cells.Item(21, 6).Value = 22
$ws1.Cells.Item(21, 7).Value = 39
$ws1.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82240'
$ws1.Cells.Item(21, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/0B72p2bL1709280677631.jpeg'

$ws1.Cells.Item(22, 2).Value = '2024-04-06'
$ws1.Cells.Item(22, 3).Value = '广州·运动番only'
$ws1.Cells.Item(22, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws1.Cells.Item(22, 5).Value = '2024.04.06 10:00-04.06 17:00'
$ws1.Cells.Item(22, 6).Value = 448
$ws1.Cells.Item(22, 7).Value = 60
$ws1.Cells.Item(22, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81454'
$ws1.Cells.Item(22, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/TBZfwnB41706255329549.jpeg'

$ws1.Cells.Item(23, 2).Value = '2024-04-06'
$ws1.Cells.Item(23, 3).Value = '广州·运动番only5.0'
$ws1.Cells.Item(23, 4).Value = '西环路1号 广州岭南会展中心'
$ws1.Cells.Item(23, 5).Value = '2024.04.06 10:00-04.06 17:00'
$ws1.Cells.Item(23, 6).Value = 40
$ws1.Cells.Item(23, 7).Value = 60
$ws1.Cells.Item(23, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82511'
$ws1.Cells.Item(23, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/JqDbCAKk1709799493319.jpeg'

$ws1.Cells.Item(24, 2).Value = '2024-04-13'
$ws1.Cells.Item(24, 3).Value = '广州·Veni Vidi Vici动漫游戏嘉年华'
$ws1.Cells.Item(24, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws1.Cells.Item(24, 5).Value = '2024.04.13 10:00-04.13 17:00'
$ws1.Cells.Item(24, 6).Value = 384
$ws1.Cells.Item(24, 7).Value = 68
$ws1.Cells.Item(24, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81575'
$ws1.Cells.Item(24, 9).Value = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'

$ws1.Cells.Item(25, 2).Value = '2024-04-13'
$ws1.Cells.Item(25, 3).Value = '广州·潮娃展WWS'
$ws1.Cells.Item(25, 4).Value = '西环路1号 广州岭南会展中心'
$ws1.Cells.Item(25, 5).Value = '2024.04.13 10:00-04.13 17:00'
$ws1.Cells.Item(25, 6).Value = 102
$ws1.Cells.Item(25, 7).Value = 48
$ws1.Cells.Item(25, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81745'
$ws1.Cells.Item(25, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/1SWNaBsA1707100228293.jpeg'

$ws1.Cells.Item(26, 2).Value = '2024-04-20'
$ws1.Cells.Item(26, 3).Value = '广州·Arknights Only·夜航星（明日方舟Only)'
$ws1.Cells.Item(26, 4).Value = '同泰路颐和山庄 颐和大酒店'
$ws1.Cells.Item(26, 5).Value = '2024.04.20 10:00-04.20 17:00'
$ws1.Cells.Item(26, 6).Value = 619
$ws1.Cells.Item(26, 7).Value = 69
$ws1.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80282'
$ws1.Cells.Item(26, 9).Value = '//i2.hdslb.com/bfs/openplatform/202312/gaEHIE1F1703745559785.jpeg'

$ws1.Cells.Item(27, 2).Value = '2024-05-04'
$ws1.Cells.Item(27, 3).Value = '广州·运动番ONLY'
$ws1.Cells.Item(27, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws1.Cells.Item(27, 5).Value = '2024.05.04 10:00-05.04 17:00'
$ws1.Cells.Item(27, 6).Value = 9
$ws1.Cells.Item(27, 7).Value = 60
$ws1.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82526'
$ws1.Cells.Item(27, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/CawIgD2O1709803813638.jpeg'

$ws1.Cells.Item(28, 2).Value = '2024-05-04'
$ws1.Cells.Item(28, 3).Value = '广州·黑塔利亚Only'
$ws1.Cells.Item(28, 4).Value = '迎宾大道123号 赛仑吉地大酒店'
$ws1.Cells.Item(28, 5).Value = '2024.05.04 09:30-05.04 16:00'
$ws1.Cells.Item(28, 6).Value = 183
$ws1.Cells.Item(28, 7).Value = 68
$ws1.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82056'
$ws1.Cells.Item(28, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/KI6tnMd81708917202487.jpeg'

$ws1.Cells.Item(29, 2).Value = '2024-05-05'
$ws1.Cells.Item(29, 3).Value = '广州·第八届萌物语动漫嘉年华'
$ws1.Cells.Item(29, 4).Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws1.Cells.Item(29, 5).Value = '2024.05.05 10:00-05.05 17:00'
$ws1.Cells.Item(29, 6).Value = 691
$ws1.Cells.Item(29, 7).Value = 60
$ws1.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81566'
$ws1.Cells.Item(29, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/c4bBhKzu1706685824726.jpeg'

$ws1.Cells.Item(30, 2).Value = '2024-05-10'
$ws1.Cells.Item(30, 3).Value = '广州·国际潮宠展—潮流创新宠物展会'
$ws1.Cells.Item(30, 4).Value = '阅江中路18号 广交会展馆C区'
$ws1.Cells.Item(30, 5).Value = '2024.05.10 10:30-05.12 18:30'
$ws1.Cells.Item(30, 6).Value = 44
$ws1.Cells.Item(30, 7).Value = 36
$ws1.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82038'
$ws1.Cells.Item(30, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/om8irfxN1708678341525.jpeg'

$ws1.Cells.Item(31, 2).Value = '2024-05-18'
$ws1.Cells.Item(31, 3).Value = '广州·恋与深空only'
$ws1.Cells.Item(31, 4).Value = '大石街石北工业大道644号 巨大创意产业园'
$ws1.Cells.Item(31, 5).Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Cells.Item(31, 6).Value = 1284
$ws1.Cells.Item(31, 7).Value = 60
$ws1.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81962'
$ws1.Cells.Item(31, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/a7aqaXrK1708485268977.jpeg'

$ws1.Cells.Item(32, 2).Value = '2024-05-18'
$ws1.Cells.Item(32, 3).Value = '广州·第五人格ONLY'
$ws1.Cells.Item(32, 4).Value = '洛浦街厦滘西环路1号 广州市岭南国际电子商务会展中心'
$ws1.Cells.Item(32, 5).Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Cells.Item(32, 6).Value = 77
$ws1.Cells.Item(32, 7).Value = 60
$ws1.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82458'
$ws1.Cells.Item(32, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/D8jK0O2X1709778592031.jpeg'

$ws1.Cells.Item(33, 2).Value = '2024-05-25'
$ws1.Cells.Item(33, 3).Value = '广州·奶司的小人国娃展Nice Mini World  '
$ws1.Cells.Item(33, 4).Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws1.Cells.Item(33, 5).Value = '2024.05.25 10:30-05.25 17:00'
$ws1.Cells.Item(33, 6).Value = 21
$ws1.Cells.Item(33, 7).Value = 60
$ws1.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82093'
$ws1.Cells.Item(33, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/rhIj7fnH1708936497981.jpeg'

# Drop the now-duplicated last row (used to be the 34th data row).
$ws1.Rows.Item(34).Delete()

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item('演出')
$ws2.Cells.Item(2, 6).Value = 356
$ws2.Cells.Item(7, 6).Value = 3

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Cells.Item(3, 6).Value = 356
$ws4.Cells.Item(8, 6).Value = 492
$ws4.Cells.Item(9, 6).Value = 492
$ws4.Cells.Item(10, 6).Value = 766
$ws4.Cells.Item(11, 6).Value = 180
$ws4.Cells.Item(12, 6).Value = 1364
$ws4.Cells.Item(13, 6).Value = 793
$ws4.Cells.Item(17, 6).Value = 597
$ws4.Cells.Item(19, 6).Value = 139
$ws4.Cells.Item(21, 6).Value = 18
$ws4.Cells.Item(23, 6).Value = 95

# Shift rows 24-46 up by one logical event, refreshing F/G counts as above.
$ws4.Cells.Item(24, 2).Value = '2024-03-30'
$ws4.Cells.Item(24, 3).Value = '广州·AP动漫游戏嘉年华'
$ws4.Cells.Item(24, 4).Value = '新港东路630-638号 南丰国际会展中心'
$ws4.Cells.Item(24, 5).Value = '2024.03.30 09:00-03.31 17:00'
$ws4.Cells.Item(24, 6).Value = 1470
$ws4.Cells.Item(24, 7).Value = '不可售'
$ws4.Cells.Item(24, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82121'
$ws4.Cells.Item(24, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/SLJ05mdG1709019165508.jpeg'

$ws4.Cells.Item(25, 2).Value = '2024-03-31'
$ws4.Cells.Item(25, 3).Value = '广州·KANAKO ITO&AYANE 2024 LIVE'
$ws4.Cells.Item(25, 4).Value = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$ws4.Cells.Item(25, 5).Value = '2024.03.31 19:00-03.31 20:30'
$ws4.Cells.Item(25, 6).Value = 192
$ws4.Cells.Item(25, 7).Value = 380
$ws4.Cells.Item(25, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81422'
$ws4.Cells.Item(25, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg'

$ws4.Cells.Item(26, 2).Value = '2024-04-04'
$ws4.Cells.Item(26, 3).Value = '广州·Look Look动漫嘉年华'
$ws4.Cells.Item(26, 4).Value = '东沙大道16号 健康方舟'
$ws4.Cells.Item(26, 5).Value = '2024.04.04 10:00-04.05 17:30'
$ws4.Cells.Item(26, 6).Value = 171
$ws4.Cells.Item(26, 7).Value = 29.9
$ws4.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82319'
$ws4.Cells.Item(26, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png'

$ws4.Cells.Item(27, 2).Value = '2024-04-05'
$ws4.Cells.Item(27, 3).Value = '广州·允月秘境动漫嘉年华'
$ws4.Cells.Item(27, 4).Value = '人民北路686号广东广播中心大楼 广东广播电视台(人民北路)'
$ws4.Cells.Item(27, 5).Value = '2024.04.05 11:30-04.06 18:00'
$ws4.Cells.Item(27, 6).Value = 22
$ws4.Cells.Item(27, 7).Value = 39
$ws4.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82240'
$ws4.Cells.Item(27, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/0B72p2bL1709280677631.jpeg'

$ws4.Cells.Item(28, 2).Value = '2024-04-06'
$ws4.Cells.Item(28, 3).Value = '广州·运动番only'
$ws4.Cells.Item(28, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws4.Cells.Item(28, 5).Value = '2024.04.06 10:00-04.06 17:00'
$ws4.Cells.Item(28, 6).Value = 448
$ws4.Cells.Item(28, 7).Value = 60
$ws4.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81454'
$ws4.Cells.Item(28, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/TBZfwnB41706255329549.jpeg'

$ws4.Cells.Item(29, 2).Value = '2024-04-06'
$ws4.Cells.Item(29, 3).Value = '广州·运动番only5.0'
$ws4.Cells.Item(29, 4).Value = '西环路1号 广州岭南会展中心'
$ws4.Cells.Item(29, 5).Value = '2024.04.06 10:00-04.06 17:00'
$ws4.Cells.Item(29, 6).Value = 40
$ws4.Cells.Item(29, 7).Value = 60
$ws4.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82511'
$ws4.Cells.Item(29, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/JqDbCAKk1709799493319.jpeg'

$ws4.Cells.Item(30, 2).Value = '2024-04-13'
$ws4.Cells.Item(30, 3).Value = '广州·Veni Vidi Vici动漫游戏嘉年华'
$ws4.Cells.Item(30, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws4.Cells.Item(30, 5).Value = '2024.04.13 10:00-04.13 17:00'
$ws4.Cells.Item(30, 6).Value = 384
$ws4.Cells.Item(30, 7).Value = 68
$ws4.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81575'
$ws4.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'

$ws4.Cells.Item(31, 2).Value = '2024-04-13'
$ws4.Cells.Item(31, 3).Value = '广州·「YOUTH NEVER GONE·直到世界尽头」演唱会'
$ws4.Cells.Item(31, 4).Value = '兴亚大道33号 广州亚运城综合体育馆'
$ws4.Cells.Item(31, 5).Value = '2024.04.13 19:00-04.13 22:00'
$ws4.Cells.Item(31, 6).Value = 3
$ws4.Cells.Item(31, 7).Value = 480
$ws4.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82607'
$ws4.Cells.Item(31, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/b1lmjEsV1709889323445.jpeg'

$ws4.Cells.Item(32, 2).Value = '2024-04-13'
$ws4.Cells.Item(32, 3).Value = '广州·潮娃展WWS'
$ws4.Cells.Item(32, 4).Value = '西环路1号 广州岭南会展中心'
$ws4.Cells.Item(32, 5).Value = '2024.04.13 10:00-04.13 17:00'
$ws4.Cells.Item(32, 6).Value = 102
$ws4.Cells.Item(32, 7).Value = 48
$ws4.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81745'
$ws4.Cells.Item(32, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/1SWNaBsA1707100228293.jpeg'

$ws4.Cells.Item(33, 2).Value = '2024-04-14'
$ws4.Cells.Item(33, 3).Value = '广州·铃木木乃美 2024 演唱会'
$ws4.Cells.Item(33, 4).Value = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$ws4.Cells.Item(33, 5).Value = '2024.04.14 19:00-04.14 20:30'
$ws4.Cells.Item(33, 6).Value = 267
$ws4.Cells.Item(33, 7).Value = 380
$ws4.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81911'
$ws4.Cells.Item(33, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/rGjpHpAV1708328728461.jpeg'

$ws4.Cells.Item(34, 2).Value = '2024-04-19'
$ws4.Cells.Item(34, 3).Value = '广州·动漫钢琴鬼才Kyle Xian互动演奏会'
$ws4.Cells.Item(34, 4).Value = '人民北路696号 广州友谊剧院'
$ws4.Cells.Item(34, 5).Value = '2024.04.19 19:30-04.19 21:00'
$ws4.Cells.Item(34, 6).Value = 61
$ws4.Cells.Item(34, 7).Value = 64
$ws4.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81434'
$ws4.Cells.Item(34, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/DJpXVLjd1706236823839.png'

$ws4.Cells.Item(35, 2).Value = '2024-04-20'
$ws4.Cells.Item(35, 3).Value = '广州·Arknights Only·夜航星（明日方舟Only)'
$ws4.Cells.Item(35, 4).Value = '同泰路颐和山庄 颐和大酒店'
$ws4.Cells.Item(35, 5).Value = '2024.04.20 10:00-04.20 17:00'
$ws4.Cells.Item(35, 6).Value = 619
$ws4.Cells.Item(35, 7).Value = 69
$ws4.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80282'
$ws4.Cells.Item(35, 9).Value = '//i2.hdslb.com/bfs/openplatform/202312/gaEHIE1F1703745559785.jpeg'

$ws4.Cells.Item(36, 2).Value = '2024-04-24'
$ws4.Cells.Item(36, 3).Value = '广州·今泉爱夏  巡演'
$ws4.Cells.Item(36, 4).Value = '革新路124号太古仓码头54汇5号仓 太空间Livehouse'
$ws4.Cells.Item(36, 5).Value = '2024.04.24 20:00-04.24 21:30'
$ws4.Cells.Item(36, 6).Value = 34
$ws4.Cells.Item(36, 7).Value = 288
$ws4.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81890'
$ws4.Cells.Item(36, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg'

$ws4.Cells.Item(37, 2).Value = '2024-04-28'
$ws4.Cells.Item(37, 3).Value = ' 广州·夏川里美 2024 巡回演唱会 出道 25 周年纪念专场'
$ws4.Cells.Item(37, 4).Value = '中山纪念堂 中山纪念堂'
$ws4.Cells.Item(37, 5).Value = '2024.04.28 19:30-04.28 21:30'
$ws4.Cells.Item(37, 6).Value = 27
$ws4.Cells.Item(37, 7).Value = 280
$ws4.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81068'
$ws4.Cells.Item(37, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/pXznRv8G1705633441713.jpeg'

$ws4.Cells.Item(38, 2).Value = '2024-04-28'
$ws4.Cells.Item(38, 3).Value = '广州·「angela LIVE 2024」in  GUANGZHOU'
$ws4.Cells.Item(38, 4).Value = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$ws4.Cells.Item(38, 5).Value = '2024.04.28 19:00-04.28 20:30'
$ws4.Cells.Item(38, 6).Value = 126
$ws4.Cells.Item(38, 7).Value = 480
$ws4.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82041'
$ws4.Cells.Item(38, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg'

$ws4.Cells.Item(39, 2).Value = '2024-04-28'
$ws4.Cells.Item(39, 3).Value = '广州·「angela LIVE 2024」in  GUANGZHOU'
$ws4.Cells.Item(39, 4).Value = '奥体南路12号优托邦购物中心 疆进酒Omni Space GZ'
$ws4.Cells.Item(39, 5).Value = '2024.04.28 19:00-04.28 20:30'
$ws4.Cells.Item(39, 6).Value = 126
$ws4.Cells.Item(39, 7).Value = 480
$ws4.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82041'
$ws4.Cells.Item(39, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg'

$ws4.Cells.Item(40, 2).Value = '2024-05-04'
$ws4.Cells.Item(40, 3).Value = '广州·运动番ONLY'
$ws4.Cells.Item(40, 4).Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws4.Cells.Item(40, 5).Value = '2024.05.04 10:00-05.04 17:00'
$ws4.Cells.Item(40, 6).Value = 9
$ws4.Cells.Item(40, 7).Value = 60
$ws4.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82526'
$ws4.Cells.Item(40, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/CawIgD2O1709803813638.jpeg'

$ws4.Cells.Item(41, 2).Value = '2024-05-04'
$ws4.Cells.Item(41, 3).Value = '广州·黑塔利亚Only'
$ws4.Cells.Item(41, 4).Value = '迎宾大道123号 赛仑吉地大酒店'
$ws4.Cells.Item(41, 5).Value = '2024.05.04 09:30-05.04 16:00'
$ws4.Cells.Item(41, 6).Value = 183
$ws4.Cells.Item(41, 7).Value = 68
$ws4.Cells.Item(41, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82056'
$ws4.Cells.Item(41, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/KI6tnMd81708917202487.jpeg'

$ws4.Cells.Item(42, 2).Value = '2024-05-05'
$ws4.Cells.Item(42, 3).Value = '广州·第八届萌物语动漫嘉年华'
$ws4.Cells.Item(42, 4).Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws4.Cells.Item(42, 5).Value = '2024.05.05 10:00-05.05 17:00'
$ws4.Cells.Item(42, 6).Value = 691
$ws4.Cells.Item(42, 7).Value = 60
$ws4.Cells.Item(42, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81566'
$ws4.Cells.Item(42, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/c4bBhKzu1706685824726.jpeg'

$ws4.Cells.Item(43, 2).Value = '2024-05-10'
$ws4.Cells.Item(43, 3).Value = '广州·国际潮宠展—潮流创新宠物展会'
$ws4.Cells.Item(43, 4).Value = '阅江中路18号 广交会展馆C区'
$ws4.Cells.Item(43, 5).Value = '2024.05.10 10:30-05.12 18:30'
$ws4.Cells.Item(43, 6).Value = 44
$ws4.Cells.Item(43, 7).Value = 36
$ws4.Cells.Item(43, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82038'
$ws4.Cells.Item(43, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/om8irfxN1708678341525.jpeg'

$ws4.Cells.Item(44, 2).Value = '2024-05-18'
$ws4.Cells.Item(44, 3).Value = '广州·恋与深空only'
$ws4.Cells.Item(44, 4).Value = '大石街石北工业大道644号 巨大创意产业园'
$ws4.Cells.Item(44, 5).Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Cells.Item(44, 6).Value = 1284
$ws4.Cells.Item(44, 7).Value = 60
$ws4.Cells.Item(44, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81962'
$ws4.Cells.Item(44, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/a7aqaXrK1708485268977.jpeg'

$ws4.Cells.Item(45, 2).Value = '2024-05-18'
$ws4.Cells.Item(45, 3).Value = '广州·第五人格ONLY'
$ws4.Cells.Item(45, 4).Value = '洛浦街厦滘西环路1号 广州市岭南国际电子商务会展中心'
$ws4.Cells.Item(45, 5).Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Cells.Item(45, 6).Value = 77
$ws4.Cells.Item(45, 7).Value = 60
$ws4.Cells.Item(45, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82458'
$ws4.Cells.Item(45, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/D8jK0O2X1709778592031.jpeg'

$ws4.Cells.Item(46, 2).Value = '2024-05-25'
$ws4.Cells.Item(46, 3).Value = '广州·奶司的小人国娃展Nice Mini World  '
$ws4.Cells.Item(46, 4).Value = '洛浦街厦滘西环路1号 岭南会展中心'
$ws4.Cells.Item(46, 5).Value = '2024.05.25 10:30-05.25 17:00'
$ws4.Cells.Item(46, 6).Value = 21
$ws4.Cells.Item(46, 7).Value = 60
$ws4.Cells.Item(46, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82093'
$ws4.Cells.Item(46, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/rhIj7fnH1708936497981.jpeg'

# Drop the now-duplicated last row (used to be the 47th data row).
$ws4.Rows.Item(47).Delete()

